$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reposicao")
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "773000"
$ws.Range("B12").Value = "REPOSIÇÃO DE PO ESPECIAL"
$ws.Range("C12").Value = "Reposicao"
